$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above row 179 (the existing row 179 and everything
# below it shift down by one, growing the used range from R292 to R293).
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row with the latest weekly price record.
$ws.Range("A179").Value = 8
$ws.Range("B179").Value = "Terminal La Palmera de La Serena"
$ws.Range("C179").Value = "Coquimbo"
$ws.Range("D179").Value = 45161
$ws.Range("E179").Value = 4
$ws.Range("F179").Value = 100112001
$ws.Range("G179").Value = "Berenjena"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 500
$ws.Range("K179").Value = 9000
$ws.Range("L179").Value = 10000
$ws.Range("M179").Value = 9500
$ws.Range("N179").Value = "`$/caja 50 unidades"
$ws.Range("O179").Value = "Región de Arica y Parinacota"
$ws.Range("P179").Value = 190
$ws.Range("Q179").Value = 50
$ws.Range("R179").Value = "Hortaliza"
